$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 2).Value = 891.2072356207559
$ws.Cells.Item(2, 3).Value = 1292.426469196052
$ws.Cells.Item(2, 4).Value = 1513.268987381566
$ws.Cells.Item(2, 5).Value = 1641.065310037129
$ws.Cells.Item(3, 2).Value = 918.8464568727627
$ws.Cells.Item(3, 3).Value = 1316.046994716525
$ws.Cells.Item(3, 4).Value = 1526.353396908353
$ws.Cells.Item(3, 5).Value = 1649.511916152012
$ws.Cells.Item(4, 2).Value = 848.1439060241668
$ws.Cells.Item(4, 3).Value = 1258.863885985375
$ws.Cells.Item(4, 4).Value = 1495.827970547305
$ws.Cells.Item(4, 5).Value = 1628.84884246994
$ws.Cells.Item(5, 2).Value = 965.1838895242677
$ws.Cells.Item(5, 3).Value = 1356.706564784675
$ws.Cells.Item(5, 4).Value = 1533.643576109059
$ws.Cells.Item(5, 5).Value = 1665.468823484424
$ws.Cells.Item(6, 2).Value = 956.0626579123599
$ws.Cells.Item(6, 3).Value = 1348.444573026424
$ws.Cells.Item(6, 4).Value = 1529.488738046546
$ws.Cells.Item(6, 5).Value = 1662.538223456871
$ws.Cells.Item(7, 2).Value = 978.3911234606398
$ws.Cells.Item(7, 3).Value = 1364.568010954196
$ws.Cells.Item(7, 4).Value = 1551.682984912907
$ws.Cells.Item(7, 5).Value = 1666.496209868257
$ws.Cells.Item(8, 2).Value = 857.551967363283
$ws.Cells.Item(8, 3).Value = 1250.574057006675
$ws.Cells.Item(8, 4).Value = 1492.356754863324
$ws.Cells.Item(8, 5).Value = 1621.016637854812
$ws.Cells.Item(9, 2).Value = 968.2246803912416
$ws.Cells.Item(9, 3).Value = 1354.716259628011
$ws.Cells.Item(9, 4).Value = 1546.03308546757
$ws.Cells.Item(9, 5).Value = 1663.019301873826
$ws.Cells.Item(10, 2).Value = 984.3029141654705
$ws.Cells.Item(10, 3).Value = 1364.70517786298
$ws.Cells.Item(10, 4).Value = 1549.14632661318
$ws.Cells.Item(10, 5).Value = 1662.741626091052
$ws.Cells.Item(11, 2).Value = 983.7696678525957
$ws.Cells.Item(11, 3).Value = 1365.677207529634
$ws.Cells.Item(11, 4).Value = 1550.616153404513
$ws.Cells.Item(11, 5).Value = 1664.139707617842
$ws.Cells.Item(12, 2).Value = 962.237953507247
$ws.Cells.Item(12, 3).Value = 1344.200000400687
$ws.Cells.Item(12, 4).Value = 1536.50930804584
$ws.Cells.Item(12, 5).Value = 1650.286811872499
$ws.Cells.Item(13, 2).Value = 982.220056285759
$ws.Cells.Item(13, 3).Value = 1362.354721856804
$ws.Cells.Item(13, 4).Value = 1547.281993772236
$ws.Cells.Item(13, 5).Value = 1660.983155509591
